$wb = $excel.ActiveWorkbook

# --- Rename sheets (drop the "20130105-" date prefix) ---
$wsDebug   = $wb.Worksheets.Item(1)
$wsRelease = $wb.Worksheets.Item(2)
$wsDebug.Name   = "PartOfSponza-Debug"
$wsRelease.Name = "PartOfSponza-Release"

# --- Update the version label shown in the Debug sheet header (C1) ---
# Shared string index 8 ("v1331") becomes unused once this cell switches
# to the already-existing "v1231" string (index 7).
$wsDebug.Range("C1").Value = "v1231"

# --- Fill in the D column of the Release sheet with the new "v1232" run ---
$wsRelease.Range("D1").Value = "v1232"
$wsRelease.Range("D2").Value  = 75
$wsRelease.Range("D3").Value  = 76
$wsRelease.Range("D4").Value  = 76
$wsRelease.Range("D5").Value  = 76
$wsRelease.Range("D6").Value  = 76
$wsRelease.Range("D7").Value  = 75
$wsRelease.Range("D8").Value  = 76
$wsRelease.Range("D9").Value  = 75
$wsRelease.Range("D10").Value = 75
$wsRelease.Range("D11").Value = 75

# Recalculate so the dependent D12:D16 summary formulas pick up real values
# instead of the previous #DIV/0! errors.
$excel.Calculate()

# --- Update the view state: active tab + selected cells ---
# The Release sheet becomes the active/selected tab.
[void]$wsDebug.Range("C2").Select()
[void]$wsRelease.Activate()
[void]$wsRelease.Range("G3").Select()
